$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Change 1: split "Per ogni piano..." paragraph run, turning
#     "4xSwitch" into "1xSwitch" and "2xAccess Point" into "1xAccess Point",
#     spread across several runs (matching the target OOXML). ---
$para2 = $d.Paragraphs(2)
$r2 = $d.Range($para2.Range.Start, $para2.Range.End)

$frag1 = '<w:p ' + $ns + '>' + `
  '<w:pPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">Per ogni piano ho preso in considerazione 30 computer massimi che si possono collegare fra loro tramite </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>1</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">xSwitch o </w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>1</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>xAccess Point e un Router per permettere l’accesso ad internet</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>.</w:t></w:r>' + `
  '</w:p>'

$r2.InsertXML($frag1)

# --- Change 2: drop the "Gli Switch..." paragraph entirely, and strip the
#     paragraph/run language formatting from the following 4 paragraphs
#     (table of subnetting + IP Network/Gateway/Broadcast). ---
$para3 = $d.Paragraphs(3)
$para7 = $d.Paragraphs(7)
$r2to7 = $d.Range($para3.Range.Start, $para7.Range.End)

$frag2 = '<w:p ' + $ns + '>' + `
  '<w:pPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>La tabella di subnetting per la rete che ho creato è la seguente:</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p ' + $ns + '><w:r><w:t>IP Network: 192.168.1.0</w:t></w:r></w:p>' + `
  '<w:p ' + $ns + '><w:r><w:t>IP Gateway: 192.168.1.1</w:t></w:r></w:p>' + `
  '<w:p ' + $ns + '><w:r><w:t>IP Broadcast: 192.168.1.255</w:t></w:r></w:p>'

$r2to7.InsertXML($frag2)
